$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.781.77"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "3.356.17"
$ws.Range("E3").Value = "  -2.02%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.07"
$ws.Range("E5").Value = "  -2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.23"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  -5.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  -4.93%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.01"
$ws.Range("E10").Value = "  -6.52%  "

$ws.Range("D11").Value = "3.353.99"
$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("E12").Value = "  -3.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.05"
$ws.Range("E13").Value = "  -3.12%  "

$ws.Range("D14").Value = "97.508.65"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.06"
$ws.Range("E15").Value = "  -5.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000254"
$ws.Range("E16").Value = "  -6.34%  "

$ws.Range("D17").Value = "3.987.22"
$ws.Range("E17").Value = "  -1.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.65"
$ws.Range("E18").Value = "  -7.36%  "

$ws.Range("D19").Value = "3.352.23"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("E21").Value = "  -14.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.85"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "509.07"
$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.34"
$ws.Range("E24").Value = "  -4.99%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.90"
$ws.Range("E25").Value = "  +7.42%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000200"
$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.25"
$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.17"
$ws.Range("E28").Value = "  -7.75%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.28"
$ws.Range("E30").Value = "  -6.53%  "

$ws.Range("E31").Value = "  -9.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.184"
$ws.Range("E32").Value = "  -7.82%  "

$ws.Range("E33").Value = "  +7.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.555"
$ws.Range("E35").Value = "  -5.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.41"
$ws.Range("E36").Value = "  -5.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.00"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "523.29"
$ws.Range("E39").Value = "  -3.20%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").Value = "  -2.71%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  +10.77%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.41"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  +2.93%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.841"
$ws.Range("E45").Value = "  -5.33%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0423"
$ws.Range("E46").Value = "  -3.92%  "

$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.61"
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.58"
$ws.Range("E48").Value = "  -7.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.46"
$ws.Range("E49").Value = "  -10.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.66"
$ws.Range("E50").Value = "  +6.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.16"
$ws.Range("E51").Value = "  -5.90%  "
